# Apply updated cryptocurrency price/volume data (and two row swaps)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.186.30'
$ws.Range("E2").Value = '  -1.03%  '
$ws.Range("D3").Value = '1.678.30'
$ws.Range("E3").Value = '  -0.94%  '
$ws.Range("E4").Value = '  -0.56%  '
$ws.Range("D5").Value = "'" + '210.88'
$ws.Range("E5").Value = '  -3.83%  '
$ws.Range("D6").Value = "'" + '0.5295'
$ws.Range("E6").Value = '  -4.55%  '
$ws.Range("D8").Value = "'" + '0.2681'
$ws.Range("E8").Value = '  -1.40%  '
$ws.Range("D9").Value = "'" + '0.06296'
$ws.Range("E9").Value = '  -2.91%  '
$ws.Range("D10").Value = "'" + '21.33'
$ws.Range("E10").Value = '  -3.77%  '
$ws.Range("D11").Value = "'" + '0.07514'
$ws.Range("E11").Value = '  -1.58%  '
$ws.Range("D12").Value = '1.671.07'
$ws.Range("E12").Value = '  -0.74%  '
$ws.Range("D13").Value = "'" + '4.486'
$ws.Range("E13").Value = '  -1.78%  '
$ws.Range("D14").Value = "'" + '0.5662'
$ws.Range("E14").Value = '  -2.90%  '
$ws.Range("D15").Value = "'" + '0.000008129'
$ws.Range("E15").Value = '  -4.01%  '
$ws.Range("D16").Value = "'" + '66.23'
$ws.Range("E16").Value = '  +1.34%  '
$ws.Range("D17").Value = '26.195.63'
$ws.Range("E17").Value = '  -1.22%  '
$ws.Range("E18").Value = '  -0.52%  '
$ws.Range("D19").Value = "'" + '4.854'
$ws.Range("E19").Value = '  -2.26%  '
$ws.Range("D20").Value = "'" + '10.53'
$ws.Range("E20").Value = '  -4.24%  '
$ws.Range("D21").Value = "'" + '188.89'
$ws.Range("E21").Value = '  -0.93%  '
$ws.Range("D22").Value = "'" + '6.199'
$ws.Range("E22").Value = '  -1.04%  '
$ws.Range("D23").Value = "'" + '1.005'
$ws.Range("E23").Value = '  -0.57%  '
$ws.Range("D24").Value = "'" + '148.13'
$ws.Range("E24").Value = '  -1.65%  '
$ws.Range("D25").Value = "'" + '0.1262'
$ws.Range("E25").Value = '  -3.84%  '
$ws.Range("D26").Value = "'" + '7.611'
$ws.Range("E26").Value = '  -3.95%  '
$ws.Range("D27").Value = "'" + '15.89'
$ws.Range("E27").Value = '  +0.75%  '
$ws.Range("D28").Value = "'" + '0.06501'
$ws.Range("E28").Value = '  +2.53%  '
$ws.Range("D29").Value = "'" + '1.342'
$ws.Range("E29").Value = '  -5.86%  '
$ws.Range("D30").Value = "'" + '1.282'
$ws.Range("E30").Value = '  -3.63%  '
$ws.Range("D31").Value = "'" + '3.527'
$ws.Range("E31").Value = '  -1.98%  '
$ws.Range("D32").Value = "'" + '3.490'
$ws.Range("E32").Value = '  -3.02%  '
$ws.Range("D33").Value = "'" + '1.646'
$ws.Range("E33").Value = '  -2.16%  '
$ws.Range("D34").Value = "'" + '1.006'
$ws.Range("E34").Value = '  -3.81%  '
$ws.Range("D35").Value = "'" + '0.6093'
$ws.Range("E35").Value = '  -2.30%  '
$ws.Range("E36").Value = '  +0.44%  '
$ws.Range("E37").Value = '  -0.44%  '
$ws.Range("D38").Value = "'" + '6.157'
$ws.Range("E38").Value = '  -1.30%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = "'" + '0.01621'
$ws.Range("E39").Value = '  -1.53%  '
$ws.Range("B40").Value = 'Maker'
$ws.Range("C40").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D40").Value = '1.100.88'
$ws.Range("E40").Value = '  -2.06%  '
$ws.Range("D41").Value = "'" + '0.8667'
$ws.Range("E41").Value = '  -2.12%  '
$ws.Range("D42").Value = "'" + '1.007'
$ws.Range("E42").Value = '  -1.03%  '
$ws.Range("D43").Value = "'" + '100.07'
$ws.Range("E43").Value = '  -0.76%  '
$ws.Range("D44").Value = '1.829.21'
$ws.Range("E44").Value = '  -0.80%  '
$ws.Range("B45").Value = 'BabyDogeCoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D45").Value = "'" + '0.00000000109'
$ws.Range("E45").Value = '  -0.86%  '
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").Value = "'" + '56.82'
$ws.Range("E46").Value = '  -1.35%  '
$ws.Range("D47").Value = "'" + '1.006'
$ws.Range("E47").Value = '  -0.21%  '
$ws.Range("D48").Value = "'" + '0.05270'
$ws.Range("E48").Value = '  -0.28%  '
$ws.Range("D49").Value = "'" + '7.979'
$ws.Range("E49").Value = '  -3.12%  '
$ws.Range("D50").Value = "'" + '0.4270'
$ws.Range("E50").Value = '  -0.76%  '
$ws.Range("D51").Value = "'" + '5.957'
$ws.Range("E51").Value = '  -2.08%  '
